$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by the
# COM Value setter (e.g. "1.00", "688.92") are written as text explicitly:
# set a temporary Text number format, assign the string, then restore the
# cell's original (unstyled/General) look via the Style property so the
# saved workbook keeps no explicit style on the cell, exactly like the rest
# of the data cells in this sheet.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '69.397.38'
$ws.Range("E2").Value = '  -2.69%  '
$ws.Range("D3").Value = '3.690.48'
$ws.Range("E3").Value = '  -3.15%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '688.92'
$ws.Range("E5").Value = '  -1.98%  '
Set-TextValue $ws.Range("D6") '162.11'
$ws.Range("E6").Value = '  -5.62%  '
$ws.Range("D7").Value = '3.688.39'
$ws.Range("E7").Value = '  -3.19%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -5.12%  '
$ws.Range("E10").Value = '  -8.68%  '
$ws.Range("E11").Value = '  -1.36%  '
$ws.Range("E12").Value = '  -6.34%  '
Set-TextValue $ws.Range("D13") '0.0000237'
$ws.Range("E13").Value = '  -6.01%  '
Set-TextValue $ws.Range("D14") '33.20'
$ws.Range("E14").Value = '  -7.71%  '
$ws.Range("D15").Value = '4.311.20'
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("D16").Value = '3.694.08'
$ws.Range("E16").Value = '  -3.84%  '
$ws.Range("D17").Value = '69.420.11'
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("E19").Value = '  -7.60%  '
Set-TextValue $ws.Range("D20") '6.54'
$ws.Range("E20").Value = '  -8.78%  '
Set-TextValue $ws.Range("D21") '477.91'
$ws.Range("E21").Value = '  -7.27%  '
Set-TextValue $ws.Range("D22") '9.97'
$ws.Range("E22").Value = '  -6.04%  '
Set-TextValue $ws.Range("D23") '0.661'
$ws.Range("E23").Value = '  -7.74%  '
Set-TextValue $ws.Range("D24") '79.93'
$ws.Range("E24").Value = '  -4.94%  '
$ws.Range("D25").Value = '3.835.88'
$ws.Range("E25").Value = '  -3.11%  '
$ws.Range("E26").Value = '  -9.88%  '
Set-TextValue $ws.Range("D27") '1.00'
$ws.Range("E27").Value = '  +0.04%  '
Set-TextValue $ws.Range("D28") '11.29'
$ws.Range("E28").Value = '  -7.03%  '
Set-TextValue $ws.Range("D29") '9.41'
$ws.Range("E29").Value = '  -8.98%  '
$ws.Range("E30").Value = '  -10.60%  '
Set-TextValue $ws.Range("D31") '2.71'
$ws.Range("E31").Value = '  -10.51%  '
Set-TextValue $ws.Range("D32") '6.81'
$ws.Range("E33").Value = '  -7.64%  '
$ws.Range("E34").Value = '  -4.93%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  -7.36%  '
$ws.Range("D37").Value = '3.661.65'
$ws.Range("E37").Value = '  -2.94%  '
Set-TextValue $ws.Range("D38") '8.42'
$ws.Range("E38").Value = '  -8.07%  '
Set-TextValue $ws.Range("D39") '6.29'
$ws.Range("E39").Value = '  +0.09%  '
Set-TextValue $ws.Range("D40") '2.33'
$ws.Range("E40").Value = '  -1.66%  '
Set-TextValue $ws.Range("D41") '0.0920'
$ws.Range("E41").Value = '  -8.87%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("E44").Value = '  -6.57%  '
Set-TextValue $ws.Range("D45") '163.72'
$ws.Range("E45").Value = '  -5.97%  '
$ws.Range("E46").Value = '  -2.68%  '
Set-TextValue $ws.Range("D47") '29.85'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("B48").Value = 'SuiNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Range("D48") '1.14'
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range("D49") '1.34'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D50") '2.75'
$ws.Range("E50").Value = '  -16.10%  '
Set-TextValue $ws.Range("D51") '0.000281'
$ws.Range("E51").Value = '  -9.33%  '
